$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.057.45"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.090.85"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.20"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.89"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "2.400.43"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.62"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.21"
$ws.Range("E14").Value = "  +3.52%  "
$ws.Range("E15").Value = "  +5.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "2.082.75"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "38.001.93"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.03"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.04"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.96"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.44"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.90"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.95"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.40"
$ws.Range("E32").Value = "  +11.06%  "
$ws.Range("E33").Value = "  +3.39%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0605"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.47"
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  +4.94%  "
$ws.Range("E38").Value = "  +7.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.13"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("D41").Value = "1.547.27"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.89"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0904"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.24"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "2.288.19"
$ws.Range("E51").Value = "  +2.77%  "
